# Entrega Final Reto 2
# Corrects the cached "Consumo de memoria" values for the "PROBING" map
# (columns B and F, rows 15-21) on the "Tablas datos" sheet, matching the
# recalculated figures from the final data run.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tablas datos")

$ws.Range("B15").Value = 107.605
$ws.Range("B16").Value = 232.732
$ws.Range("B17").Value = 502.91
$ws.Range("B18").Value = 798.71100000000001
$ws.Range("B19").Value = 1400.3030000000001
$ws.Range("B20").Value = 2435.2829999999999
$ws.Range("B21").Value = 3297.1689999999999

$ws.Range("F16").Value = 239.696
$ws.Range("F17").Value = 493.37299999999999
$ws.Range("F18").Value = 795.67
$ws.Range("F19").Value = 1407.2070000000001
$ws.Range("F20").Value = 2400.0329999999999
$ws.Range("F21").Value = 3112.1680000000001
